$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting existing columns (B..AI) to (C..AJ).
$ws.Columns("B").Insert()

# New header + data for the inserted "gstNumber" column.
$ws.Range("B1").Value = "gstNumber"
$ws.Range("B2").Value = "17CMRPS9572E1Z7"
$ws.Range("B3").Value = "17CMRPS9572E1Z7"
$ws.Range("B4").Value = "17CMRPS9572E1Z7"

# Page setup (paper size / orientation) matching the saved workbook's print settings.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the selection Excel leaves on save.
$ws.Range("B5").Select()
